$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()
Write-Host "Inserted"
